$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1349
$ws.Range("I2").Value = 3266
$ws.Range("J2").Value = 13928
$ws.Range("K2").Value = 66
$ws.Range("L2").Value = 3871
$ws.Range("M2").Value = 243
$ws.Range("N2").Value = 2472
$ws.Range("O2").Value = 9
$ws.Range("P2").Value = 51
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = 183
$ws.Range("S2").Value = 1513
$ws.Range("T2").Value = 2488
$ws.Range("U2").Value = 180
$ws.Range("V2").Value = 21797
$ws.Range("X2").Value = 21757
$ws.Range("Z2").Value = 327
$ws.Range("AA2").Value = 130
